# arreglado el error del saldo Null
# Fill in the previously-empty "saldo" (column H) cells with 0 so the
# balance column no longer shows a Null/blank value for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locales")

$rows = @(3, 6, 7, 9, 12, 18, 19, 24, 25, 27, 29, 30, 37, 39, 40, 42, 43, 44, 45, 46, 50, 52, 53, 54, 55, 60, 62)

foreach ($r in $rows) {
    $ws.Range("H$r").Value = 0
}

# Reflect where the editor ended up after making the fix: scrolled further
# down the sheet with the active cell resting on H64.
$ws.Activate()
$ws.Range("F34").Select()
$ws.Range("H64").Select()
